# Apply the "NL Stats" rename + updated session data to stats.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the data on both worksheets (they currently hold identical data).
#    Row 5 (Xavier) is removed; rows 2-4 get new numbers; "Cedric" becomes
#    "Raymond" with a new date in T3.
# ---------------------------------------------------------------------------
$sheetNames = @("combined Stats-this session", "PLO Stats-this session")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2 (Fish)
    $ws.Range("C2").Value = 53.59
    $ws.Range("D2").Value = 33.59
    $ws.Range("F2").Value = 0.679
    $ws.Range("G2").Value = 0
    $ws.Range("I2").Value = 0.118
    $ws.Range("J2").Value = 0.268
    $ws.Range("K2").Value = 0.19
    $ws.Range("L2").Value = 0.19
    $ws.Range("N2").Value = 0
    $ws.Range("O2").Value = 112.17
    $ws.Range("P2").Value = 66.58
    $ws.Range("Q2").Value = 168
    $ws.Range("R2").Value = 0.711

    # Row 3 (Cedric -> Raymond)
    $ws.Range("A3").Value = "Raymond"
    $ws.Range("B3").Value = 51
    $ws.Range("C3").Value = 24.39
    $ws.Range("D3").Value = -26.61
    $ws.Range("F3").Value = 0.466
    $ws.Range("G3").Value = 0.31
    $ws.Range("H3").Value = 0.0057
    $ws.Range("I3").Value = 0.454
    $ws.Range("J3").Value = 0.126
    $ws.Range("K3").Value = 0.046
    $ws.Range("L3").Value = 3.68
    $ws.Range("M3").Value = 26
    $ws.Range("N3").Value = 56
    $ws.Range("O3").Value = 41.7
    $ws.Range("P3").Value = 58.54
    $ws.Range("Q3").Value = 174
    $ws.Range("R3").Value = 0.364
    # Force text (not an auto-converted date serial) to match the original
    # inline-string cell style, like typing it with a leading apostrophe.
    $ws.Range("T3").Value = "'07/05/21"

    # Row 4 (Scott)
    $ws.Range("C4").Value = 13.02
    $ws.Range("D4").Value = -6.98
    $ws.Range("F4").Value = 0.527
    $ws.Range("G4").Value = 0.24
    $ws.Range("H4").Value = 0.02
    $ws.Range("I4").Value = 0.513
    $ws.Range("J4").Value = 0.22
    $ws.Range("K4").Value = 0.08
    $ws.Range("L4").Value = 3.16
    $ws.Range("M4").Value = 16
    $ws.Range("N4").Value = 45
    $ws.Range("O4").Value = 43.65
    $ws.Range("P4").Value = 64.11
    $ws.Range("Q4").Value = 150
    $ws.Range("R4").Value = 0.364

    # Row 5 (Xavier) is dropped entirely from the session.
    $ws.Rows("5:5").Delete()
}

# ---------------------------------------------------------------------------
# 2) Fix up every chart's series formulas so they still point at the right
#    (now 4-row) ranges. Deleting the row above does not retarget the chart
#    series automatically, so each one is rewritten explicitly.
# ---------------------------------------------------------------------------
for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $chartObjects = $ws.ChartObjects()
    for ($i = 1; $i -le $chartObjects.Count; $i++) {
        $chart = $chartObjects.Item($i).Chart
        $series = $chart.SeriesCollection()
        for ($j = 1; $j -le $series.Count; $j++) {
            $ser = $series.Item($j)
            $formula = $ser.Formula
            $formula = $formula.Replace("PLO Stats-this session", "NL Stats-this session")
            $formula = $formula.Replace('$5', '$4')
            $ser.Formula = $formula
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Rename the second tab from "PLO Stats-this session" to
#    "NL Stats-this session" (do this last so the lookups above by the old
#    name still resolve).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("PLO Stats-this session").Name = "NL Stats-this session"
